# Update the participant-count / category-count figures in Table 2.3
# (inversion of ML and Activity Count fix, per commit message).
#
# Each entry is (1-based table row, old cell text, new cell text) for the
# single 2-column table in the document.  Row 17 ("(372,586]") keeps its
# original "22,564 (25)" value and is intentionally NOT in this list, even
# though the same text also appears (and IS changed) on rows 13 and 14.
#
# Note: cells are addressed by row/column and their Range.Text is set
# directly (rather than via Find/Replace) so that rows sharing identical
# text (e.g. "22,564 (25)" appears on rows 13, 14 and 17) are updated
# independently and precisely, without affecting sibling cells that
# happen to contain the same string.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1;  Old = "N = 90,258";   New = "N = 90,237" },
    @{ Row = 3;  Old = "89,943 (100)"; New = "89,922 (100)" },
    @{ Row = 4;  Old = "89,612 (99)";  New = "89,592 (99)" },
    @{ Row = 6;  Old = "22,594 (25)";  New = "22,589 (25)" },
    @{ Row = 7;  Old = "22,588 (25)";  New = "22,584 (25)" },
    @{ Row = 8;  Old = "22,518 (25)";  New = "22,516 (25)" },
    @{ Row = 9;  Old = "22,558 (25)";  New = "22,548 (25)" },
    @{ Row = 11; Old = "22,576 (25)";  New = "22,573 (25)" },
    @{ Row = 12; Old = "22,554 (25)";  New = "22,550 (25)" },
    @{ Row = 13; Old = "22,564 (25)";  New = "22,557 (25)" },
    @{ Row = 14; Old = "22,564 (25)";  New = "22,557 (25)" },
    @{ Row = 16; Old = "22,566 (25)";  New = "22,560 (25)" },
    @{ Row = 18; Old = "22,563 (25)";  New = "22,558 (25)" },
    @{ Row = 19; Old = "22,565 (25)";  New = "22,559 (25)" },
    @{ Row = 23; Old = "2,017 (2.2)";  New = "2,016 (2.2)" },
    @{ Row = 24; Old = "1,284 (1.4)";  New = "1,283 (1.4)" },
    @{ Row = 27; Old = "51,733 (57)";  New = "51,723 (57)" },
    @{ Row = 28; Old = "38,525 (43)";  New = "38,514 (43)" },
    @{ Row = 30; Old = "7,120 (7.9)";  New = "7,118 (7.9)" },
    @{ Row = 31; Old = "22,117 (25)";  New = "22,114 (25)" },
    @{ Row = 32; Old = "21,321 (24)";  New = "21,318 (24)" },
    @{ Row = 33; Old = "39,700 (44)";  New = "39,687 (44)" }
)

foreach ($change in $changes) {
    $cell = $t.Cell($change.Row, 2)
    $rng = $cell.Range
    # Cell.Range.Text carries the trailing cell-mark (CR + cell-end char);
    # trim it before comparing against the plain expected text.
    $current = $rng.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $change.Old) {
        Write-Host "WARNING: row" $change.Row "expected" $change.Old "but found" $current
    }
    $rng.Text = $change.New
}
